$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

# Columns A (Date) and D (Week) contain values that look like a date / number
# ("2023-06-09", "23"). Force them to be stored as text (matching the rest of
# the column), then drop the temporary number-format style so the new row
# keeps the workbook's default cell style, same as every other data row.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-09"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "09:02:34"
$ws.Cells.Item($row, 3).Value = "Friday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "23"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 120087
$ws.Cells.Item($row, 6).Value = 134589
$ws.Cells.Item($row, 7).Value = 160242
$ws.Cells.Item($row, 8).Value = 131501
$ws.Cells.Item($row, 9).Value = 175597
$ws.Cells.Item($row, 10).Value = 113591
$ws.Cells.Item($row, 11).Value = 201204
$ws.Cells.Item($row, 12).Value = 221226
$ws.Cells.Item($row, 13).Value = 173181
$ws.Cells.Item($row, 14).Value = 120271
$ws.Cells.Item($row, 15).Value = 38615
$ws.Cells.Item($row, 16).Value = 34410
$ws.Cells.Item($row, 17).Value = 50883
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36361
$ws.Cells.Item($row, 20).Value = -1
